# Add multi-browser support and custom annotations for test categorization
# and authorship; implement retry logic for failed tests.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("TestsRunner")

# --- TestsRunner sheet: new "Browser" column, new "username"/"password"
# --- annotation columns, and a retry row for loginLogoutTest ---

# Header row
$ws2.Range("A1").Value = "TestCase"
$ws2.Range("B1").Value = "Description"
$ws2.Range("C1").Value = "Browser"
$ws2.Range("D1").Value = "Execute"
$ws2.Range("E1").Value = "username"
$ws2.Range("F1").Value = "password"

# loginLogoutTest - chrome run
$ws2.Range("A2").Value = "loginLogoutTest"
$ws2.Range("B2").Value = "validate OrangeHRM login and logout functionality"
$ws2.Range("C2").Value = "chrome"
$ws2.Range("D2").Value = "no"
$ws2.Range("E2").Value = "admin"
$ws2.Range("F2").Value = "ad123"

# loginLogoutTest - retry row (chrome)
$ws2.Range("A3").Value = "loginLogoutTest"
$ws2.Range("B3").Value = "validate OrangeHRM login and logout functionality"
$ws2.Range("C3").Value = "chrome"
$ws2.Range("D3").Value = "no"
$ws2.Range("E3").Value = "ad123"
$ws2.Range("F3").Value = "admin"

# loginLogoutTest - firefox run
$ws2.Range("A4").Value = "loginLogoutTest"
$ws2.Range("B4").Value = "validate OrangeHRM login and logout functionality"
$ws2.Range("C4").Value = "firefox"
$ws2.Range("D4").Value = "yes"
$ws2.Range("E4").Value = "Admin"
$ws2.Range("F4").Value = "admin123"

# homePageTitleTest - chrome run
$ws2.Range("A5").Value = "homePageTitleTest"
$ws2.Range("B5").Value = "validate title of home page"
$ws2.Range("C5").Value = "chrome"
$ws2.Range("D5").Value = "yes"
$ws2.Range("E5").Value = "Admin"
$ws2.Range("F5").Value = "admin123"

# Widen the new Browser column
$ws2.Columns.Item(3).ColumnWidth = 15.3

# --- Selections ---
$ws1.Activate() | Out-Null
$ws1.Range("C8").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("D11").Select() | Out-Null
